$d = $word.ActiveDocument
for ($i = 174; $i -le 187; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    Write-Output ($i.ToString() + " len=" + $t.Length + " text=[" + $t + "]")
}
